$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap column widths: column A takes the former column B width, and vice versa
$ws.Columns.Item(1).ColumnWidth = 15.7109375
$ws.Columns.Item(2).ColumnWidth = 16.42578125

# Update cell values per diff (rows 1-32, columns A and B)
$ws.Range("A1").Value = -0.27204752559465817
$ws.Range("B1").Value = 0.27132573398009185
$ws.Range("A2").Value = -0.23655335530744903
$ws.Range("B2").Value = 0.23376104813802367
$ws.Range("A3").Value = -0.11215656413347652
$ws.Range("B3").Value = 0.11174184728325898
$ws.Range("A4").Value = -0.10374184731934122
$ws.Range("B4").Value = 0.10337470744447508
$ws.Range("A5").Value = -0.10037470746621224
$ws.Range("B5").Value = 0.099137416619453411
$ws.Range("A6").Value = 0.00015926418011069643
$ws.Range("B6").Value = -0.00018722022942618821
$ws.Range("A7").Value = 0.010187220180425616
$ws.Range("B7").Value = -0.010189242609685234
$ws.Range("A8").Value = 0.020189242561128307
$ws.Range("B8").Value = -0.020219342338518853
$ws.Range("A9").Value = 0.022219342317611801
$ws.Range("B9").Value = -0.022263873225557784
$ws.Range("A10").Value = 0.024263873206221476
$ws.Range("B10").Value = -0.024268890639364216
$ws.Range("A11").Value = 0.027268890616809038
$ws.Range("B11").Value = -0.027282396720417879
$ws.Range("A12").Value = -0.020861923775481461
$ws.Range("B12").Value = 0.020669585067941476
$ws.Range("A13").Value = -0.017169585093185447
$ws.Range("B13").Value = 0.017081854729313228
$ws.Range("A14").Value = -0.0090818547697288921
$ws.Range("B14").Value = 0.0090530663337275996
$ws.Range("A15").Value = -0.0080530663514295497
$ws.Range("B15").Value = 0.0080344637740159541
$ws.Range("A16").Value = -0.0060344637952645108
$ws.Range("B16").Value = 0.0060035890943601622
$ws.Range("A17").Value = -0.004003589115957773
$ws.Range("B17").Value = 0.0039999999717243995
$ws.Range("A18").Value = -0.0287723786935139
$ws.Range("B18").Value = 0.028734492540355205
$ws.Range("A19").Value = -0.012092118693931564
$ws.Range("B19").Value = 0.012017023446798447
$ws.Range("A20").Value = -0.0080170234621608216
$ws.Range("B20").Value = 0.0080057083895184178
$ws.Range("A21").Value = -0.0040057084050477698
$ws.Range("B21").Value = 0.0039999999843525202
$ws.Range("A22").Value = -0.11760448402537804
$ws.Range("B22").Value = 0.11664208800202402
$ws.Range("A23").Value = -0.090448713038565387
$ws.Range("B23").Value = 0.089332539033468272
$ws.Range("A24").Value = -0.020102540805486235
$ws.Range("B24").Value = 0.019999999917858169
$ws.Range("A25").Value = -0.097296680822610782
$ws.Range("B25").Value = 0.097168901145700204
$ws.Range("A26").Value = -0.094668901171122144
$ws.Range("B26").Value = 0.094504731412799359
$ws.Range("A27").Value = -0.092004731439771614
$ws.Range("B27").Value = 0.091034908390279057
$ws.Range("A28").Value = -0.089034908421890435
$ws.Range("B28").Value = 0.088370822877638311
$ws.Range("A29").Value = -0.081370822930185049
$ws.Range("B29").Value = 0.081178086086198142
$ws.Range("A30").Value = -0.021178086314936007
$ws.Range("B30").Value = 0.021024913210321738
$ws.Range("A31").Value = -0.014024913266720063
$ws.Range("B31").Value = 0.014001499389820538
$ws.Range("A32").Value = -0.004001499456268931
$ws.Range("B32").Value = 0.0039999999531623587
